$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ncam1"
$ws.Cells.Item(2,3).Value = "Gfra1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.7004376666666666
$ws.Cells.Item(2,8).Value = 2.101313
$ws.Cells.Item(2,9).Value = 0.04511966030063898
$ws.Cells.Item(2,10).Value = 0.04511966030063898
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.06547366666666667
$ws.Cells.Item(2,14).Value = 0.196421
$ws.Cells.Item(2,15).Value = 0.002125877360986814
$ws.Cells.Item(2,16).Value = 0.002125877360986814
$ws.Cells.Item(2,17).Value = 0.04586022230811111
$ws.Cells.Item(2,18).Value = 0.412742000773
$ws.Cells.Item(2,19).Value = 0.00009591886436854393
$ws.Cells.Item(2,20).Value = 0.00009591886436854392

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ncam1"
$ws.Cells.Item(3,3).Value = "Gfra1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.7004376666666666
$ws.Cells.Item(3,8).Value = 2.101313
$ws.Cells.Item(3,9).Value = 0.04511966030063898
$ws.Cells.Item(3,10).Value = 0.04511966030063898
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 25.94643066666667
$ws.Cells.Item(3,14).Value = 77.839292
$ws.Cells.Item(3,15).Value = 0.8424597607080814
$ws.Cells.Item(3,16).Value = 0.8424597607080814
$ws.Cells.Item(3,17).Value = 18.17385735448844
$ws.Cells.Item(3,18).Value = 163.564716190396
$ws.Cells.Item(3,19).Value = 0.03801149822010624
$ws.Cells.Item(3,20).Value = 0.03801149822010623

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ncam1"
$ws.Cells.Item(4,3).Value = "Gfra1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.7004376666666666
$ws.Cells.Item(4,8).Value = 2.101313
$ws.Cells.Item(4,9).Value = 0.04511966030063898
$ws.Cells.Item(4,10).Value = 0.04511966030063898
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 4.786517
$ws.Cells.Item(4,14).Value = 14.359551
$ws.Cells.Item(4,15).Value = 0.1554143619309319
$ws.Cells.Item(4,16).Value = 0.1554143619309319
$ws.Cells.Item(4,17).Value = 3.352656798940333
$ws.Cells.Item(4,18).Value = 30.173911190463
$ws.Cells.Item(4,19).Value = 0.007012243216164205
$ws.Cells.Item(4,20).Value = 0.007012243216164204

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ncam1"
$ws.Cells.Item(5,3).Value = "Gfra1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.237305
$ws.Cells.Item(5,8).Value = 0.711915
$ws.Cells.Item(5,9).Value = 0.01528632952964618
$ws.Cells.Item(5,10).Value = 0.01528632952964618
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.06547366666666667
$ws.Cells.Item(5,14).Value = 0.196421
$ws.Cells.Item(5,15).Value = 0.002125877360986814
$ws.Cells.Item(5,16).Value = 0.002125877360986814
$ws.Cells.Item(5,17).Value = 0.01553722846833333
$ws.Cells.Item(5,18).Value = 0.139835056215
$ws.Cells.Item(5,19).Value = 0.00003249686187965903
$ws.Cells.Item(5,20).Value = 0.00003249686187965903

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ncam1"
$ws.Cells.Item(6,3).Value = "Gfra1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.237305
$ws.Cells.Item(6,8).Value = 0.711915
$ws.Cells.Item(6,9).Value = 0.01528632952964618
$ws.Cells.Item(6,10).Value = 0.01528632952964618
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 25.94643066666667
$ws.Cells.Item(6,14).Value = 77.839292
$ws.Cells.Item(6,15).Value = 0.8424597607080814
$ws.Cells.Item(6,16).Value = 0.8424597607080814
$ws.Cells.Item(6,17).Value = 6.157217729353333
$ws.Cells.Item(6,18).Value = 55.41495956417999
$ws.Cells.Item(6,19).Value = 0.0128781175176506
$ws.Cells.Item(6,20).Value = 0.0128781175176506

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ncam1"
$ws.Cells.Item(7,3).Value = "Gfra1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.237305
$ws.Cells.Item(7,8).Value = 0.711915
$ws.Cells.Item(7,9).Value = 0.01528632952964618
$ws.Cells.Item(7,10).Value = 0.01528632952964618
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 4.786517
$ws.Cells.Item(7,14).Value = 14.359551
$ws.Cells.Item(7,15).Value = 0.1554143619309319
$ws.Cells.Item(7,16).Value = 0.1554143619309319
$ws.Cells.Item(7,17).Value = 1.135864416685
$ws.Cells.Item(7,18).Value = 10.222779750165
$ws.Cells.Item(7,19).Value = 0.002375715150115923
$ws.Cells.Item(7,20).Value = 0.002375715150115923

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Ncam1"
$ws.Cells.Item(8,3).Value = "Gfra1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 14.58625866666667
$ws.Cells.Item(8,8).Value = 43.758776
$ws.Cells.Item(8,9).Value = 0.9395940101697148
$ws.Cells.Item(8,10).Value = 0.9395940101697148
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.06547366666666667
$ws.Cells.Item(8,14).Value = 0.196421
$ws.Cells.Item(8,15).Value = 0.002125877360986814
$ws.Cells.Item(8,16).Value = 0.002125877360986814
$ws.Cells.Item(8,17).Value = 0.955015837855111
$ws.Cells.Item(8,18).Value = 8.595142540696001
$ws.Cells.Item(8,19).Value = 0.001997461634738611
$ws.Cells.Item(8,20).Value = 0.001997461634738611

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Ncam1"
$ws.Cells.Item(9,3).Value = "Gfra1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 14.58625866666667
$ws.Cells.Item(9,8).Value = 43.758776
$ws.Cells.Item(9,9).Value = 0.9395940101697148
$ws.Cells.Item(9,10).Value = 0.9395940101697148
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 25.94643066666667
$ws.Cells.Item(9,14).Value = 77.839292
$ws.Cells.Item(9,15).Value = 0.8424597607080814
$ws.Cells.Item(9,16).Value = 0.8424597607080814
$ws.Cells.Item(9,17).Value = 378.4613491807324
$ws.Cells.Item(9,18).Value = 3406.152142626592
$ws.Cells.Item(9,19).Value = 0.7915701449703245
$ws.Cells.Item(9,20).Value = 0.7915701449703245

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Ncam1"
$ws.Cells.Item(10,3).Value = "Gfra1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 14.58625866666667
$ws.Cells.Item(10,8).Value = 43.758776
$ws.Cells.Item(10,9).Value = 0.9395940101697148
$ws.Cells.Item(10,10).Value = 0.9395940101697148
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 4.786517
$ws.Cells.Item(10,14).Value = 14.359551
$ws.Cells.Item(10,15).Value = 0.1554143619309319
$ws.Cells.Item(10,16).Value = 0.1554143619309319
$ws.Cells.Item(10,17).Value = 69.81737507439733
$ws.Cells.Item(10,18).Value = 628.356375669576
$ws.Cells.Item(10,19).Value = 0.1460264035646517
$ws.Cells.Item(10,20).Value = 0.1460264035646517
